$p = $ppt.ActivePresentation
$newDate = "13.08.2018"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -like "*2018*") {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# All slide layouts
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}

# Notes master (direct shape text writes are not persisted for the
# notes master in this host; go through HeadersFooters.DateAndTime instead)
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = $newDate
